# Adds new columns I ("I0") and J ("IF") to the worksheet, matching the
# header style already used by column H, and fills in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (font/border/alignment) from the existing "IP"
# header (H1) onto the two new header cells so they match the other
# headers in row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Data values (rows 2-73) ------------------------------------------
# Each entry is @(I-value, J-value) for the corresponding data row,
# starting at row 2.
$data = @(
    @(6,6),
    @(8,8),
    @(6,6),
    @(7,7),
    @(5,5),
    @(9,9),
    @(8,8),
    @(8,8),
    @(6,6),
    @(8,8),
    @(7,7),
    @(6,6),
    @(5,5),
    @(10,10),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(6,7),
    @(11,11),
    @(6,6),
    @(10,10),
    @(7,7),
    @(5,5),
    @(7,7),
    @(8,8),
    @(6,6),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(2,2),
    @(7,7),
    @(9,9),
    @(7,8),
    @(7,7),
    @(5,6),
    @(4,4),
    @(6,6),
    @(6,7),
    @(7,7),
    @(8,8),
    @(6,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(9,9),
    @(6,6),
    @(9,9),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(6,7),
    @(7,7),
    @(9,9),
    @(7,7),
    @(7,7),
    @(9,9),
    @(7,7),
    @(7,7),
    @(7,7),
    @(9,9),
    @(9,9),
    @(9,9),
    @(4,4),
    @(5,5),
    @(6,6),
    @(5,5),
    @(6,6),
    @(5,5),
    @(3,3)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
